$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1420.2
$ws.Range("J6").Value = 2166.6667
$ws.Range("L6").Value = 6500.000100000001
$ws.Range("N6").Value = -6724.000100000001
$ws.Range("H8").Value = 19.166666
$ws.Range("I8").Value = 19.166666
$ws.Range("K8").Value = 57.499998
$ws.Range("M8").Value = 81.50000199999999
$ws.Range("H32").Value = 11835.111
$ws.Range("I32").Value = 15507
$ws.Range("J32").Value = 9999.166999999999
$ws.Range("K32").Value = 15507
$ws.Range("L32").Value = 9999.166999999999
$ws.Range("M32").Value = -15181
$ws.Range("N32").Value = -10651.167
$ws.Range("H125").Value = 1759
$ws.Range("I125").Value = 1800
$ws.Range("K125").Value = 16200
$ws.Range("M125").Value = -13740
$ws.Range("H131").Value = 4910.5
$ws.Range("J131").Value = 17263
$ws.Range("L131").Value = 51789
$ws.Range("N131").Value = -61869
$ws.Range("H132").Value = 1542.697
$ws.Range("I132").Value = 1372.0714
$ws.Range("K132").Value = 4116.2142
$ws.Range("M132").Value = -1586.2142
$ws.Range("H137").Value = 3846.238
$ws.Range("I137").Value = 4033.3333
$ws.Range("K137").Value = 12099.9999
$ws.Range("M137").Value = -9549.999899999999
$ws.Range("H138").Value = 3690.2292
$ws.Range("I138").Value = 3130.5715
$ws.Range("J138").Value = 3920.6765
$ws.Range("K138").Value = 9391.7145
$ws.Range("L138").Value = 11762.0295
$ws.Range("M138").Value = -4251.7145
$ws.Range("N138").Value = -22042.0295
$ws.Range("H141").Value = 856
$ws.Range("I141").Value = 856
$ws.Range("K141").Value = 2568
$ws.Range("M141").Value = 2612

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1551.9
$ws.Range("I32").Value = 1196.5526
$ws.Range("K32").Value = 1196.5526
$ws.Range("M32").Value = -909.5526
$ws.Range("H45").Value = 66669604
$ws.Range("I45").Value = 83334480
$ws.Range("J45").Value = 10102.667
$ws.Range("K45").Value = 83334480
$ws.Range("L45").Value = 10102.667
$ws.Range("M45").Value = -83334103
$ws.Range("N45").Value = -10856.667
$ws.Range("H61").Value = 6169.1797
$ws.Range("I61").Value = 5069.567
$ws.Range("J61").Value = 9834.556
$ws.Range("K61").Value = 5069.567
$ws.Range("L61").Value = 9834.556
$ws.Range("M61").Value = -4857.567
$ws.Range("N61").Value = -10258.556
$ws.Range("H74").Value = 10418168
$ws.Range("I74").Value = 13890104
$ws.Range("K74").Value = 13890104
$ws.Range("M74").Value = -13889230
$ws.Range("H77").Value = 10418168
$ws.Range("I77").Value = 13890104
$ws.Range("K77").Value = 69450520
$ws.Range("M77").Value = -69446152
$ws.Range("H102").Value = 1385.0476
$ws.Range("I102").Value = 1449.2222
$ws.Range("K102").Value = 1449.2222
$ws.Range("M102").Value = 172.7778000000001
$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680
$ws.Range("H111").Value = 73492.664
$ws.Range("J111").Value = 73492.664
$ws.Range("L111").Value = 73492.664
$ws.Range("N111").Value = -81672.664
$ws.Range("H114").Value = 69994.5
$ws.Range("J114").Value = 69994.5
$ws.Range("L114").Value = 69994.5
$ws.Range("N114").Value = -78672.5
$ws.Range("H136").Value = 6169.1797
$ws.Range("I136").Value = 5069.567
$ws.Range("J136").Value = 9834.556
$ws.Range("K136").Value = 15208.701
$ws.Range("L136").Value = 29503.668
$ws.Range("M136").Value = -12658.701
$ws.Range("N136").Value = -34603.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9456.9375
$ws.Range("I105").Value = 4631.2
$ws.Range("J105").Value = 17499.834
$ws.Range("K105").Value = 4631.2
$ws.Range("L105").Value = 17499.834
$ws.Range("M105").Value = -2884.2
$ws.Range("N105").Value = -20993.834
$ws.Range("H107").Value = 997.1429000000001
$ws.Range("I107").Value = 496.83334
$ws.Range("K107").Value = 496.83334
$ws.Range("M107").Value = 1423.16666
$ws.Range("H134").Value = 3546.3333
$ws.Range("I134").Value = 1712.7142
$ws.Range("J134").Value = 9964
$ws.Range("K134").Value = 5138.142599999999
$ws.Range("L134").Value = 29892
$ws.Range("M134").Value = -2603.142599999999
$ws.Range("N134").Value = -34962
$ws.Range("H135").Value = 61999.2
$ws.Range("J135").Value = 61999.2
$ws.Range("L135").Value = 61999.2
$ws.Range("N135").Value = -72139.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38978.234
$ws.Range("I31").Value = 3348.238
$ws.Range("K31").Value = 3348.238
$ws.Range("M31").Value = -3053.238
$ws.Range("H34").Value = 38978.234
$ws.Range("I34").Value = 3348.238
$ws.Range("K34").Value = 3348.238
$ws.Range("M34").Value = -3146.238
$ws.Range("H104").Value = 49995
$ws.Range("J104").Value = 49995
$ws.Range("L104").Value = 49995
$ws.Range("N104").Value = -55237
$ws.Range("H132").Value = 4473.0527
$ws.Range("I132").Value = 4361.7144
$ws.Range("J132").Value = 4784.8
$ws.Range("K132").Value = 13085.1432
$ws.Range("L132").Value = 14354.4
$ws.Range("M132").Value = -10555.1432
$ws.Range("N132").Value = -19414.4
$ws.Range("H134").Value = 2035.6765
$ws.Range("I134").Value = 1333.625
$ws.Range("K134").Value = 4000.875
$ws.Range("M134").Value = -1465.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1985.4546
$ws.Range("I25").Value = 2295
$ws.Range("J25").Value = 1916.6666
$ws.Range("K25").Value = 6885
$ws.Range("L25").Value = 5749.9998
$ws.Range("M25").Value = -6716
$ws.Range("N25").Value = -6087.9998
$ws.Range("H30").Value = 1985.4546
$ws.Range("I30").Value = 2295
$ws.Range("J30").Value = 1916.6666
$ws.Range("K30").Value = 6885
$ws.Range("L30").Value = 5749.9998
$ws.Range("M30").Value = -6783
$ws.Range("N30").Value = -5953.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2108.6667
$ws.Range("I102").Value = 2108.6667
$ws.Range("K102").Value = 2108.6667
$ws.Range("M102").Value = -486.6667000000002
$ws.Range("H114").Value = 38500
$ws.Range("J114").Value = 38500
$ws.Range("L114").Value = 38500
$ws.Range("N114").Value = -47178
$ws.Range("H132").Value = 4810.9116
$ws.Range("I132").Value = 2905.7036
$ws.Range("K132").Value = 8717.110799999999
$ws.Range("M132").Value = -6187.110799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2641.75
$ws.Range("I16").Value = 2376.2856
$ws.Range("J16").Value = 4500
$ws.Range("K16").Value = 2376.2856
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = -2206.2856
$ws.Range("N16").Value = -4840
$ws.Range("H22").Value = 2681.625
$ws.Range("I22").Value = 1145.4
$ws.Range("J22").Value = 5242
$ws.Range("K22").Value = 1145.4
$ws.Range("L22").Value = 5242
$ws.Range("M22").Value = -850.4000000000001
$ws.Range("N22").Value = -5832
$ws.Range("H27").Value = 2681.625
$ws.Range("I27").Value = 1145.4
$ws.Range("J27").Value = 5242
$ws.Range("K27").Value = 1145.4
$ws.Range("L27").Value = 5242
$ws.Range("M27").Value = -1038.4
$ws.Range("N27").Value = -5456
$ws.Range("H40").Value = 8337.1875
$ws.Range("J40").Value = 12127.5
$ws.Range("L40").Value = 12127.5
$ws.Range("N40").Value = -12399.5
$ws.Range("H68").Value = 5349
$ws.Range("I68").Value = 1935.5
$ws.Range("K68").Value = 1935.5
$ws.Range("M68").Value = -1186.5
$ws.Range("H71").Value = 5349
$ws.Range("I71").Value = 1935.5
$ws.Range("K71").Value = 9677.5
$ws.Range("M71").Value = -5933.5
$ws.Range("H82").Value = 6656.143
$ws.Range("J82").Value = 5548.125
$ws.Range("L82").Value = 5548.125
$ws.Range("N82").Value = -6270.125
$ws.Range("H85").Value = 6656.143
$ws.Range("J85").Value = 5548.125
$ws.Range("L85").Value = 5548.125
$ws.Range("N85").Value = -8044.125
$ws.Range("H136").Value = 10501
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 27534.5
$ws.Range("J47").Value = 27534.5
$ws.Range("L47").Value = 27534.5
$ws.Range("N47").Value = -28678.5
$ws.Range("H49").Value = 19531
$ws.Range("J49").Value = 19531
$ws.Range("L49").Value = 19531
$ws.Range("N49").Value = -19991
$ws.Range("H54").Value = 24519.25
$ws.Range("J54").Value = 24519.25
$ws.Range("L54").Value = 24519.25
$ws.Range("N54").Value = -25559.25
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 35000
$ws.Range("L70").Value = 35000
$ws.Range("N70").Value = -35630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("N73").Value = -37184
$ws.Range("H106").Value = 55000
$ws.Range("J106").Value = 55000
$ws.Range("L106").Value = 55000
$ws.Range("N106").Value = -57524
$ws.Range("H113").Value = 923
$ws.Range("I113").Value = 699.5
$ws.Range("J113").Value = 1072
$ws.Range("K113").Value = 2098.5
$ws.Range("L113").Value = 3216
$ws.Range("M113").Value = 71.5
$ws.Range("N113").Value = -7556
